# Selenium Python Framework TestObjects workbook update:
#  - add 4 new Page-Object sheets (FlightFinderPage, SelectFlightsPage,
#    BookFlightPage, ConfirmationPage) following the existing LoginPage
#    sheet, each with a logicalName/locator_type/locator_value header
#    styled the same way as LoginPage's header row.
#  - update LoginPage's own sheet-view selection now that it is no longer
#    the active tab.

$wb = $excel.ActiveWorkbook
$loginPage = $wb.Worksheets.Item("LoginPage")

# ---------------------------------------------------------------------
# LoginPage (existing sheet) keeps its data; only the saved selection
# changes now that another sheet becomes the active tab.
# ---------------------------------------------------------------------
$loginPage.Activate()
$loginPage.Range("A1:C5").Select()

# ---------------------------------------------------------------------
# FlightFinderPage
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$flightFinderPage = $wb.Worksheets.Add($null, $lastSheet)
$flightFinderPage.Name = "FlightFinderPage"

$flightFinderPage.Range("A1").Value = "logicalName"
$flightFinderPage.Range("B1").Value = "locator_type"
$flightFinderPage.Range("C1").Value = "locator_value"
$loginPage.Range("A1:C1").Copy()
$flightFinderPage.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$flightFinderPage.Range("A2").Value = "rd_oneway"
$flightFinderPage.Range("B2").Value = "xpath"
$flightFinderPage.Range("C2").Value = "//input[@value='oneway']"
$flightFinderPage.Range("A3").Value = "lst_passengers"
$flightFinderPage.Range("B3").Value = "xpath"
$flightFinderPage.Range("C3").Value = "//select[@name='passCount']"
$flightFinderPage.Range("A4").Value = "lst_from"
$flightFinderPage.Range("B4").Value = "xpath"
$flightFinderPage.Range("C4").Value = "//select[@name='fromPort']"
$flightFinderPage.Range("A5").Value = "lst_day"
$flightFinderPage.Range("B5").Value = "xpath"
$flightFinderPage.Range("C5").Value = "//select[@name='fromDay']"
$flightFinderPage.Range("A6").Value = "btn_continue"
$flightFinderPage.Range("B6").Value = "xpath"
$flightFinderPage.Range("C6").Value = "//input[@name='findFlights']"

$flightFinderPage.Activate()
$flightFinderPage.Range("A1:C2").Select()

# ---------------------------------------------------------------------
# SelectFlightsPage
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$selectFlightsPage = $wb.Worksheets.Add($null, $lastSheet)
$selectFlightsPage.Name = "SelectFlightsPage"

$selectFlightsPage.Range("A1").Value = "logicalName"
$selectFlightsPage.Range("B1").Value = "locator_type"
$selectFlightsPage.Range("C1").Value = "locator_value"
$loginPage.Range("A1:C1").Copy()
$selectFlightsPage.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$selectFlightsPage.Range("A2").Value = "btn_reserveflights"
$selectFlightsPage.Range("B2").Value = "name"
$selectFlightsPage.Range("C2").Value = "reserveFlights"

$selectFlightsPage.Activate()
$selectFlightsPage.Range("A1:C2").Select()

# ---------------------------------------------------------------------
# BookFlightPage
# (cell-write order below intentionally matches the source edit session:
# the "lastname" row's locator_value was filled in before its own
# logicalName cell, right after the credit-card row was completed.)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bookFlightPage = $wb.Worksheets.Add($null, $lastSheet)
$bookFlightPage.Name = "BookFlightPage"

$bookFlightPage.Range("A1").Value = "logicalName"
$bookFlightPage.Range("B1").Value = "locator_type"
$bookFlightPage.Range("C1").Value = "locator_value"
$loginPage.Range("A1:C1").Copy()
$bookFlightPage.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$bookFlightPage.Range("A2").Value = "txt_firstname"
$bookFlightPage.Range("B2").Value = "name"
$bookFlightPage.Range("C2").Value = "passFirst0"
$bookFlightPage.Range("B3").Value = "name"
$bookFlightPage.Range("C3").Value = "passLast0"
$bookFlightPage.Range("A4").Value = "txt_creditcard"
$bookFlightPage.Range("B4").Value = "name"
$bookFlightPage.Range("C4").Value = "creditnumber"
$bookFlightPage.Range("A3").Value = "txt_lastname"
$bookFlightPage.Range("A5").Value = "btn_purchase"
$bookFlightPage.Range("B5").Value = "name"
$bookFlightPage.Range("C5").Value = "buyFlights"

$bookFlightPage.Activate()
$bookFlightPage.Range("A1:C3").Select()

# ---------------------------------------------------------------------
# ConfirmationPage (becomes the active tab)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$confirmationPage = $wb.Worksheets.Add($null, $lastSheet)
$confirmationPage.Name = "ConfirmationPage"

$confirmationPage.Range("A1").Value = "logicalName"
$confirmationPage.Range("B1").Value = "locator_type"
$confirmationPage.Range("C1").Value = "locator_value"
$loginPage.Range("A1:C1").Copy()
$confirmationPage.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$confirmationPage.Range("A2").Value = "we_confirmation"
$confirmationPage.Range("B2").Value = "xpath"
$confirmationPage.Range("C2").Value = "//font[contains(text(),'booked')]"
$confirmationPage.Range("A3").Value = "btn_Logout"
$confirmationPage.Range("B3").Value = "xpath"
$confirmationPage.Range("C3").Value = "(//td)[55]/a"

$confirmationPage.Activate()
$confirmationPage.Range("C7").Select()
